$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 2 and row 3: A, B, E, F, G, I, M
# The edit swaps these values between the two rows.
# Note: use Value2 for reading - Value has a bug in this runtime that
# stringifies the underlying .NET object instead of returning the value.
$plainCols = @("A", "B", "E", "F", "G", "M")

foreach ($col in $plainCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value = $val3
    $cell3.Value = $val2
}

# Column I holds numeric-looking text ("1"/"2") that must stay text
# (not get auto-converted to a real number). Force text formatting,
# assign, then clear the formatting so no stray number-format style
# is left behind on the cell while the stored type remains text.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = [string]$i3
$ws.Range("I2").ClearFormats()

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = [string]$i2
$ws.Range("I3").ClearFormats()
